$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '59.285.92'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.579.53'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '573.53'
$ws.Range('E5').Value = '  +3.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.33'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.597'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.588.05'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  +2.92%  '
$ws.Range('E12').Value = '  +11.39%  '
$ws.Range('E13').Value = '  +3.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.032.10'
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '59.298.20'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '22.50'
$ws.Range('E16').Value = '  +7.28%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000138'
$ws.Range('E17').Value = '  +3.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.585.09'
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '339.34'
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.27'
$ws.Range('E21').Value = '  +1.69%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.27'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '64.76'
$ws.Range('E24').Value = '  -2.69%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.460'
$ws.Range('E25').Value = '  +7.33%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.27'
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0784'
$ws.Range('E29').Value = '  +3.20%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.69'
$ws.Range('E31').Value = '  +0.71%  '
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '158.12'
$ws.Range('E33').Value = '  +2.54%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '19.04'
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.04'
$ws.Range('E35').Value = '  +1.85%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.16'
$ws.Range('E36').Value = '  +2.08%  '
$ws.Range('E37').Value = '  -3.78%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.878'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '37.24'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '295.99'
$ws.Range('E41').Value = '  +4.51%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.67'
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0977'
$ws.Range('E44').Value = '  +2.42%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '128.18'
$ws.Range('E45').Value = '  +8.21%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.594'
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0537'
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '19.23'
$ws.Range('E48').Value = '  +2.35%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '10.63'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('E50').Value = '  +2.71%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.954.22'
$ws.Range('E51').Value = '  +0.32%  '
